# ---------------------------------------------------------------------------
# CS133JS_Lab07_CodeReview.docx edit script
#
# 1. Table 2 (Web App I), row 12 ("... file called as a result of clicking
#    the button?"): split the trailing run so "as a result of" is wrapped in
#    gramStart/gramEnd proofErr markers, then insert a brand-new row right
#    after it asking "Is code that does i/o separated from code that does
#    processing?" (with the "i" wrapped in spellStart/spellEnd markers).
# 2. Table 3 (Web App II) heading row: collapse the many runs that spell out
#    "Web App II: Multiple-choice quiz, true-false quiz, fill-in-the-blank
#    quiz" into a single run.
# 3. Table 3, row 12 (same question as #1): apply the same gramStart/gramEnd
#    split (no new row here, it was already added once in table 2... wait,
#    table 3 also needs its own i/o row? No - only table 2 gets a new row).
# 4. Table 3, row 13 ("Is document.querySelector used at least once?"):
#    remove the stray _GoBack bookmark and merge "document." + "querySelector"
#    into a single run, moving the gramStart/gramEnd markers so they wrap the
#    whole "document.querySelector" phrase (nested inside spellStart/End).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$W14 = "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

# --- Change 1: Table 2, row 12 -------------------------------------------
$t2 = $d.Tables.Item(2)
$cell = $t2.Rows.Item(12).Cells.Item(1)
$xml = "<w:p $W $W14 w14:paraId='74608713' w14:textId='5BC0ADF6' w:rsidR='00F10BB5' w:rsidRDefault='00A26FEB' w:rsidP='00E50705'>" + `
  "<w:pPr><w:tabs><w:tab w:val='left' w:pos='0'/></w:tabs><w:suppressAutoHyphens w:val='0'/><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:rPr><w:bCs/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space='preserve'>   Is a function in the .</w:t></w:r>" + `
  "<w:proofErr w:type='spellStart'/><w:r><w:rPr><w:bCs/></w:rPr><w:t>js</w:t></w:r><w:proofErr w:type='spellEnd'/>" + `
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space='preserve'> file called </w:t></w:r>" + `
  "<w:proofErr w:type='gramStart'/><w:r><w:rPr><w:bCs/></w:rPr><w:t>as a result of</w:t></w:r><w:proofErr w:type='gramEnd'/>" + `
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space='preserve'> clicking the button?</w:t></w:r>" + `
  "</w:p>"
$cell.Range.InsertXML($xml)

# Insert a brand-new row right after row 12 for the i/o separation question.
$newRow = $t2.Rows.Add($t2.Rows.Item(13))
$xmlNew = "<w:p $W><w:pPr><w:tabs><w:tab w:val='left' w:pos='0'/></w:tabs><w:suppressAutoHyphens w:val='0'/><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:rPr><w:bCs/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space='preserve'>   Is code that does </w:t></w:r>" + `
  "<w:proofErr w:type='spellStart'/><w:r><w:rPr><w:bCs/></w:rPr><w:t>i</w:t></w:r><w:proofErr w:type='spellEnd'/>" + `
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t>/o separated from code that does processing?</w:t></w:r>" + `
  "</w:p>"
$newRow.Cells.Item(1).Range.InsertXML($xmlNew)

# --- Change 2: Table 3, row 1 - collapse "Web App II..." heading runs -----
$t3 = $d.Tables.Item(3)
$headingCell = $t3.Rows.Item(1).Cells.Item(1)
$xmlHeading = "<w:p $W $W14 w14:paraId='42D60F46' w14:textId='33EAF653' w:rsidR='00777666' w:rsidRPr='002A1596' w:rsidRDefault='00777666' w:rsidP='00E06635'>" + `
  "<w:pPr><w:widowControl w:val='0'/><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:rPr><w:b/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:b/></w:rPr><w:t>Web App II: Multiple-choice quiz, true-false quiz, fill-in-the-blank quiz</w:t></w:r>" + `
  "</w:p>"
$headingCell.Range.InsertXML($xmlHeading)

# --- Change 3: Table 3, row 12 - same gramStart/gramEnd split -------------
$cell3 = $t3.Rows.Item(12).Cells.Item(1)
$xml3 = "<w:p $W $W14 w14:paraId='72C9C8B4' w14:textId='77777777' w:rsidR='00777666' w:rsidRDefault='00777666' w:rsidP='00E06635'>" + `
  "<w:pPr><w:tabs><w:tab w:val='left' w:pos='0'/></w:tabs><w:suppressAutoHyphens w:val='0'/><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:rPr><w:bCs/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space='preserve'>   Is a function in the .</w:t></w:r>" + `
  "<w:proofErr w:type='spellStart'/><w:r><w:rPr><w:bCs/></w:rPr><w:t>js</w:t></w:r><w:proofErr w:type='spellEnd'/>" + `
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space='preserve'> file called </w:t></w:r>" + `
  "<w:proofErr w:type='gramStart'/><w:r><w:rPr><w:bCs/></w:rPr><w:t>as a result of</w:t></w:r><w:proofErr w:type='gramEnd'/>" + `
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space='preserve'> clicking the button?</w:t></w:r>" + `
  "</w:p>"
$cell3.Range.InsertXML($xml3)

# --- Change 4: Table 3, row 13 - drop _GoBack bookmark, merge runs --------
$cell4 = $t3.Rows.Item(13).Cells.Item(1)
$xml4 = "<w:p $W $W14 w14:paraId='1E832AB7' w14:textId='625EBF93' w:rsidR='00777666' w:rsidRDefault='00777666' w:rsidP='00E06635'>" + `
  "<w:pPr><w:tabs><w:tab w:val='left' w:pos='0'/></w:tabs><w:suppressAutoHyphens w:val='0'/><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:rPr><w:bCs/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space='preserve'>   Is </w:t></w:r>" + `
  "<w:proofErr w:type='spellStart'/><w:proofErr w:type='gramStart'/>" + `
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t>document.querySelector</w:t></w:r>" + `
  "<w:proofErr w:type='spellEnd'/><w:proofErr w:type='gramEnd'/>" + `
  "<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space='preserve'> used at least once?</w:t></w:r>" + `
  "</w:p>"
$cell4.Range.InsertXML($xml4)
